# Update existing rows 2-9 and add new rows 10-13 to reflect the re-run
# of the natmi LR-pair analysis (Il1b-Il1r1) with an extra replicate
# ("sCs" sending cluster) folded in, per Dr Hou's advice.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Il1b"
$row2[0,2] = "Il1r1"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 0.197811
$row2[0,7] = 0.593433
$row2[0,8] = 0.001274141111268169
$row2[0,9] = 0.001274141111268169
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 13.21972033333333
$row2[0,13] = 39.659161
$row2[0,14] = 0.2311669015805739
$row2[0,15] = 0.2311669015805739
$row2[0,16] = 2.615006098857
$row2[0,17] = 23.535054889713
$row2[0,18] = 0.0002945392528682919
$row2[0,19] = 0.0002945392528682919
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Il1b"
$row3[0,2] = "Il1r1"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 0.197811
$row3[0,7] = 0.593433
$row3[0,8] = 0.001274141111268169
$row3[0,9] = 0.001274141111268169
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 34.595189
$row3[0,13] = 103.785567
$row3[0,14] = 0.6049494580123129
$row3[0,15] = 0.6049494580123129
$row3[0,16] = 6.843308931278999
$row3[0,17] = 61.589780381511
$row3[0,18] = 0.0007707909746928848
$row3[0,19] = 0.0007707909746928849
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Il1b"
$row4[0,2] = "Il1r1"
$row4[0,3] = "M2"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 0.197811
$row4[0,7] = 0.593433
$row4[0,8] = 0.001274141111268169
$row4[0,9] = 0.001274141111268169
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 0.110909
$row4[0,13] = 0.332727
$row4[0,14] = 0.001939412426354648
$row4[0,15] = 0.001939412426354648
$row4[0,16] = 0.021939020199
$row4[0,17] = 0.197451181791
$row4[0,18] = 0.000002471085104122806
$row4[0,19] = 0.000002471085104122807
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Il1b"
$row5[0,2] = "Il1r1"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 0.197811
$row5[0,7] = 0.593433
$row5[0,8] = 0.001274141111268169
$row5[0,9] = 0.001274141111268169
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 9.261089666666665
$row5[0,13] = 27.783269
$row5[0,14] = 0.1619442279807586
$row5[0,15] = 0.1619442279807586
$row5[0,16] = 1.831945408053
$row5[0,17] = 16.487508672477
$row5[0,18] = 0.0002063397986028694
$row5[0,19] = 0.0002063397986028695
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "M2"
$row6[0,1] = "Il1b"
$row6[0,2] = "Il1r1"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 154.0066273333333
$row6[0,7] = 462.0198820000001
$row6[0,8] = 0.9919881871744044
$row6[0,9] = 0.9919881871744045
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 13.21972033333333
$row6[0,13] = 39.659161
$row6[0,14] = 0.2311669015805739
$row6[0,15] = 0.2311669015805739
$row6[0,16] = 2035.924542826556
$row6[0,17] = 18323.320885439
$row6[0,18] = 0.2293148356336375
$row6[0,19] = 0.2293148356336375
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "M2"
$row7[0,1] = "Il1b"
$row7[0,2] = "Il1r1"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 154.0066273333333
$row7[0,7] = 462.0198820000001
$row7[0,8] = 0.9919881871744044
$row7[0,9] = 0.9919881871744045
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 34.595189
$row7[0,13] = 103.785567
$row7[0,14] = 0.6049494580123129
$row7[0,15] = 0.6049494580123129
$row7[0,16] = 5327.888379849232
$row7[0,17] = 47950.9954186431
$row7[0,18] = 0.6001027161857728
$row7[0,19] = 0.6001027161857728
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "M2"
$row8[0,1] = "Il1b"
$row8[0,2] = "Il1r1"
$row8[0,3] = "M2"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 154.0066273333333
$row8[0,7] = 462.0198820000001
$row8[0,8] = 0.9919881871744044
$row8[0,9] = 0.9919881871744045
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 0.110909
$row8[0,13] = 0.332727
$row8[0,14] = 0.001939412426354648
$row8[0,15] = 0.001939412426354648
$row8[0,16] = 17.08072103091267
$row8[0,17] = 153.726489278214
$row8[0,18] = 0.00192387421700306
$row8[0,19] = 0.00192387421700306
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "M2"
$row9[0,1] = "Il1b"
$row9[0,2] = "Il1r1"
$row9[0,3] = "sCs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 154.0066273333333
$row9[0,7] = 462.0198820000001
$row9[0,8] = 0.9919881871744044
$row9[0,9] = 0.9919881871744045
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 9.261089666666665
$row9[0,13] = 27.783269
$row9[0,14] = 0.1619442279807586
$row9[0,15] = 0.1619442279807586
$row9[0,16] = 1426.269184994917
$row9[0,17] = 12836.42266495426
$row9[0,18] = 0.1606467611379911
$row9[0,19] = 0.1606467611379912
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "sCs"
$row10[0,1] = "Il1b"
$row10[0,2] = "Il1r1"
$row10[0,3] = "ECs"
$row10[0,4] = 1
$row10[0,5] = 0.3333333333333333
$row10[0,6] = 1.046026666666667
$row10[0,7] = 3.13808
$row10[0,8] = 0.006737671714327339
$row10[0,9] = 0.006737671714327339
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 13.21972033333333
$row10[0,13] = 39.659161
$row10[0,14] = 0.2311669015805739
$row10[0,15] = 0.2311669015805739
$row10[0,16] = 13.82817999454222
$row10[0,17] = 124.45361995088
$row10[0,18] = 0.001557526694068125
$row10[0,19] = 0.001557526694068124
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "sCs"
$row11[0,1] = "Il1b"
$row11[0,2] = "Il1r1"
$row11[0,3] = "FAPs"
$row11[0,4] = 1
$row11[0,5] = 0.3333333333333333
$row11[0,6] = 1.046026666666667
$row11[0,7] = 3.13808
$row11[0,8] = 0.006737671714327339
$row11[0,9] = 0.006737671714327339
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 34.595189
$row11[0,13] = 103.785567
$row11[0,14] = 0.6049494580123129
$row11[0,15] = 0.6049494580123129
$row11[0,16] = 36.18749023237333
$row11[0,17] = 325.68741209136
$row11[0,18] = 0.004075950851847215
$row11[0,19] = 0.004075950851847215
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "sCs"
$row12[0,1] = "Il1b"
$row12[0,2] = "Il1r1"
$row12[0,3] = "M2"
$row12[0,4] = 1
$row12[0,5] = 0.3333333333333333
$row12[0,6] = 1.046026666666667
$row12[0,7] = 3.13808
$row12[0,8] = 0.006737671714327339
$row12[0,9] = 0.006737671714327339
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 0.110909
$row12[0,13] = 0.332727
$row12[0,14] = 0.001939412426354648
$row12[0,15] = 0.001939412426354648
$row12[0,16] = 0.1160137715733333
$row12[0,17] = 1.04412394416
$row12[0,18] = 0.00001306712424746466
$row12[0,19] = 0.00001306712424746466
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "sCs"
$row13[0,1] = "Il1b"
$row13[0,2] = "Il1r1"
$row13[0,3] = "sCs"
$row13[0,4] = 1
$row13[0,5] = 0.3333333333333333
$row13[0,6] = 1.046026666666667
$row13[0,7] = 3.13808
$row13[0,8] = 0.006737671714327339
$row13[0,9] = 0.006737671714327339
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 9.261089666666665
$row13[0,13] = 27.783269
$row13[0,14] = 0.1619442279807586
$row13[0,15] = 0.1619442279807586
$row13[0,16] = 9.687346753724443
$row13[0,17] = 87.18612078352
$row13[0,18] = 0.001091127044164535
$row13[0,19] = 0.001091127044164535
$ws.Range("A13:T13").Value = $row13

